$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 384, shifting existing rows 384-401 down to 389-406
$ws.Range("A384:R388").Insert()

# Row 384
$ws.Cells.Item(384,1).Value = 6
$ws.Cells.Item(384,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(384,3).Value = "Metropolitana"
$ws.Cells.Item(384,4).Value = 44509
$ws.Cells.Item(384,5).Value = 13
$ws.Cells.Item(384,6).Value = 100112003
$ws.Cells.Item(384,7).Value = "Ajo"
$ws.Cells.Item(384,8).Value = "Chino"
$ws.Cells.Item(384,9).Value = "1a nueva(o)"
$ws.Cells.Item(384,10).Value = 35000
$ws.Cells.Item(384,11).Value = 2000
$ws.Cells.Item(384,12).Value = 2000
$ws.Cells.Item(384,13).Value = 2000
$ws.Cells.Item(384,14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(384,15).Value = "Provincia de Talagante"
$ws.Cells.Item(384,16).Value = 100
$ws.Cells.Item(384,17).Value = 20
$ws.Cells.Item(384,18).Value = "Hortaliza"

# Row 385
$ws.Cells.Item(385,1).Value = 6
$ws.Cells.Item(385,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(385,3).Value = "Metropolitana"
$ws.Cells.Item(385,4).Value = 44509
$ws.Cells.Item(385,5).Value = 13
$ws.Cells.Item(385,6).Value = 100112003
$ws.Cells.Item(385,7).Value = "Ajo"
$ws.Cells.Item(385,8).Value = "Chino"
$ws.Cells.Item(385,9).Value = "2a nueva(o)"
$ws.Cells.Item(385,10).Value = 27000
$ws.Cells.Item(385,11).Value = 1600
$ws.Cells.Item(385,12).Value = 1600
$ws.Cells.Item(385,13).Value = 1600
$ws.Cells.Item(385,14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(385,15).Value = "Provincia de Talagante"
$ws.Cells.Item(385,16).Value = 80
$ws.Cells.Item(385,17).Value = 20
$ws.Cells.Item(385,18).Value = "Hortaliza"

# Row 386
$ws.Cells.Item(386,1).Value = 6
$ws.Cells.Item(386,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(386,3).Value = "Metropolitana"
$ws.Cells.Item(386,4).Value = 44509
$ws.Cells.Item(386,5).Value = 13
$ws.Cells.Item(386,6).Value = 100112003
$ws.Cells.Item(386,7).Value = "Ajo"
$ws.Cells.Item(386,8).Value = "Chino"
$ws.Cells.Item(386,9).Value = "3a nueva (o)"
$ws.Cells.Item(386,10).Value = 15000
$ws.Cells.Item(386,11).Value = 800
$ws.Cells.Item(386,12).Value = 800
$ws.Cells.Item(386,13).Value = 800
$ws.Cells.Item(386,14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(386,15).Value = "Provincia de Talagante"
$ws.Cells.Item(386,16).Value = 40
$ws.Cells.Item(386,17).Value = 20
$ws.Cells.Item(386,18).Value = "Hortaliza"

# Row 387
$ws.Cells.Item(387,1).Value = 6
$ws.Cells.Item(387,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(387,3).Value = "Metropolitana"
$ws.Cells.Item(387,4).Value = 44509
$ws.Cells.Item(387,5).Value = 13
$ws.Cells.Item(387,6).Value = 100112003
$ws.Cells.Item(387,7).Value = "Ajo"
$ws.Cells.Item(387,8).Value = "Chino"
$ws.Cells.Item(387,9).Value = "Extra nueva (o)"
$ws.Cells.Item(387,10).Value = 31000
$ws.Cells.Item(387,11).Value = 2500
$ws.Cells.Item(387,12).Value = 2500
$ws.Cells.Item(387,13).Value = 2500
$ws.Cells.Item(387,14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(387,15).Value = "Provincia de Talagante"
$ws.Cells.Item(387,16).Value = 125
$ws.Cells.Item(387,17).Value = 20
$ws.Cells.Item(387,18).Value = "Hortaliza"

# Row 388
$ws.Cells.Item(388,1).Value = 6
$ws.Cells.Item(388,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(388,3).Value = "Metropolitana"
$ws.Cells.Item(388,4).Value = 44509
$ws.Cells.Item(388,5).Value = 13
$ws.Cells.Item(388,6).Value = 100112003
$ws.Cells.Item(388,7).Value = "Ajo"
$ws.Cells.Item(388,8).Value = "Chino"
$ws.Cells.Item(388,9).Value = "Primera"
$ws.Cells.Item(388,10).Value = 1700
$ws.Cells.Item(388,11).Value = 16500
$ws.Cells.Item(388,12).Value = 17000
$ws.Cells.Item(388,13).Value = 16824
$ws.Cells.Item(388,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(388,15).Value = "China"
$ws.Cells.Item(388,16).Value = 1682
$ws.Cells.Item(388,17).Value = 10
$ws.Cells.Item(388,18).Value = "Hortaliza"
